$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8606349999999999
$ws.Range("H2").Value = 2.581905
$ws.Range("I2").Value = 0.0262626340301864
$ws.Range("J2").Value = 0.0262626340301864
$ws.Range("M2").Value = 4.093680666666667
$ws.Range("N2").Value = 12.281042
$ws.Range("O2").Value = 0.1610908176055751
$ws.Range("P2").Value = 0.161090817605575
$ws.Range("Q2").Value = 3.523164860556666
$ws.Range("R2").Value = 31.70848374501
$ws.Range("S2").Value = 0.004230669188398726
$ws.Range("T2").Value = 0.004230669188398725
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8606349999999999
$ws.Range("H3").Value = 2.581905
$ws.Range("I3").Value = 0.0262626340301864
$ws.Range("J3").Value = 0.0262626340301864
$ws.Range("O3").Value = 0.5606512265211691
$ws.Range("P3").Value = 0.5606512265211691
$ws.Range("Q3").Value = 12.261820565985
$ws.Range("R3").Value = 110.356385093865
$ws.Range("S3").Value = 0.0147241779807006
$ws.Range("T3").Value = 0.0147241779807006
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8606349999999999
$ws.Range("H4").Value = 2.581905
$ws.Range("I4").Value = 0.0262626340301864
$ws.Range("J4").Value = 0.0262626340301864
$ws.Range("M4").Value = 7.071161666666666
$ws.Range("N4").Value = 21.213485
$ws.Range("O4").Value = 0.2782579558732559
$ws.Range("P4").Value = 0.2782579558732559
$ws.Range("Q4").Value = 6.085689220991665
$ws.Range("R4").Value = 54.771202988925
$ws.Range("S4").Value = 0.007307786861087074
$ws.Range("T4").Value = 0.007307786861087073
$ws.Range("H5").Value = 58.40949000000001
$ws.Range("I5").Value = 0.5941299388474139
$ws.Range("J5").Value = 0.5941299388474139
$ws.Range("M5").Value = 4.093680666666667
$ws.Range("N5").Value = 12.281042
$ws.Range("O5").Value = 0.1610908176055751
$ws.Range("P5").Value = 0.161090817605575
$ws.Range("Q5").Value = 79.70326665428668
$ws.Range("R5").Value = 717.32939988858
$ws.Range("S5").Value = 0.09570887761288022
$ws.Range("T5").Value = 0.0957088776128802
$ws.Range("H6").Value = 58.40949000000001
$ws.Range("I6").Value = 0.5941299388474139
$ws.Range("J6").Value = 0.5941299388474139
$ws.Range("O6").Value = 0.5606512265211691
$ws.Range("P6").Value = 0.5606512265211691
$ws.Range("R6").Value = 2496.55203099117
$ws.Range("S6").Value = 0.3330996789277498
$ws.Range("T6").Value = 0.3330996789277498
$ws.Range("H7").Value = 58.40949000000001
$ws.Range("I7").Value = 0.5941299388474139
$ws.Range("J7").Value = 0.5941299388474139
$ws.Range("M7").Value = 7.071161666666666
$ws.Range("N7").Value = 21.213485
$ws.Range("O7").Value = 0.2782579558732559
$ws.Range("P7").Value = 0.2782579558732559
$ws.Range("Q7").Value = 137.6743155525167
$ws.Range("R7").Value = 1239.06883997265
$ws.Range("S7").Value = 0.1653213823067839
$ws.Range("T7").Value = 0.1653213823067839
$ws.Range("G8").Value = 12.439858
$ws.Range("H8").Value = 37.319574
$ws.Range("I8").Value = 0.3796074271223998
$ws.Range("J8").Value = 0.3796074271223997
$ws.Range("M8").Value = 4.093680666666667
$ws.Range("N8").Value = 12.281042
$ws.Range("O8").Value = 0.1610908176055751
$ws.Range("P8").Value = 0.161090817605575
$ws.Range("Q8").Value = 50.92480619067867
$ws.Range("R8").Value = 458.323255716108
$ws.Range("S8").Value = 0.06115127080429612
$ws.Range("T8").Value = 0.06115127080429611
$ws.Range("G9").Value = 12.439858
$ws.Range("H9").Value = 37.319574
$ws.Range("I9").Value = 0.3796074271223998
$ws.Range("J9").Value = 0.3796074271223997
$ws.Range("O9").Value = 0.5606512265211691
$ws.Range("P9").Value = 0.5606512265211691
$ws.Range("Q9").Value = 177.235769707638
$ws.Range("R9").Value = 1595.121927368742
$ws.Range("S9").Value = 0.2128273696127187
$ws.Range("T9").Value = 0.2128273696127187
$ws.Range("G10").Value = 12.439858
$ws.Range("H10").Value = 37.319574
$ws.Range("I10").Value = 0.3796074271223998
$ws.Range("J10").Value = 0.3796074271223997
$ws.Range("M10").Value = 7.071161666666666
$ws.Range("N10").Value = 21.213485
$ws.Range("O10").Value = 0.2782579558732559
$ws.Range("P10").Value = 0.2782579558732559
$ws.Range("Q10").Value = 87.96424702837665
$ws.Range("R10").Value = 791.6782232553899
$ws.Range("S10").Value = 0.1056287867053849
$ws.Range("T10").Value = 0.1056287867053849
